$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = "#00af8c"
$ws.Range("B18").Value = "#8a73b4"
$ws.Range("B32").Value = "#8C8D91"
$ws.Range("C32").Value = "square"
$ws.Range("B33").Value = "#00ABF5"
$ws.Range("C33").Value = "diamond"
$ws.Range("B34").Value = "#0AB6FF"
$ws.Range("B35").Value = "#33C2FF"
$ws.Range("B36").Value = "#7C134B"
$ws.Range("B37").Value = "#A01860"
$ws.Range("B38").Value = "#00E0E0"
$ws.Range("B39").Value = "#5CCEFF"
$ws.Range("B40").Value = "#A8A9AC"
$ws.Range("B41").Value = "#E75FA8"
$ws.Range("B42").Value = "#B0B1B4"
$ws.Range("B43").Value = "#B8B9BC"
$ws.Range("B44").Value = "#E972B1"
$ws.Range("B45").Value = "#EC83BB"
$ws.Range("B46").Value = "#73D039"
$ws.Range("B47").Value = "#98999D"
$ws.Range("B48").Value = "#88898D"
$ws.Range("B49").Value = "#A0A1A4"
$ws.Range("B50").Value = "#C0C1C4"
$ws.Range("B51").Value = "#4C4345"
$ws.Range("B52").Value = "#DA3036"
$ws.Range("B53").Value = "#008FCC"
$ws.Range("B54").Value = "#D8D4D4"
$ws.Range("C54").Value = "circle"
$ws.Range("B55").Value = "#084427"
$ws.Range("C55").Value = "circle"
$ws.Range("B56").Value = "#006635"
$ws.Range("C56").Value = "circle"
$ws.Range("B57").Value = "#008F4A"
$ws.Range("C57").Value = "circle"
$ws.Range("B58").Value = "#00CC69"
$ws.Range("C58").Value = "circle"
$ws.Range("B59").Value = "#48706B"
$ws.Range("C59").Value = "diamond"
$ws.Range("B60").Value = "#009B9B"
$ws.Range("C60").Value = "square"
$ws.Range("B61").Value = "#E4595D"
$ws.Range("C61").Value = "diamond"
$ws.Range("B62").Value = "#00D6D6"
$ws.Range("C62").Value = "square"
$ws.Range("B63").Value = "#47FFFF"
$ws.Range("B64").Value = "#9C9DA1"
$ws.Range("C64").Value = "square"
$ws.Range("B65").Value = "#85DAFF"
$ws.Range("C65").Value = "diamond"
$ws.Range("B66").Value = "#0072A3"
$ws.Range("C66").Value = "diamond"
$ws.Range("B67").Value = "#00F5F5"
$ws.Range("B68").Value = "#47C8FF"
$ws.Range("C68").Value = "diamond"
$ws.Range("B69").Value = "#00af8c"
$ws.Range("B70").Value = "#413A3B"
$ws.Range("B71").Value = "#78797D"
$ws.Range("B72").Value = "#E23C94"
$ws.Range("B73").Value = "#008F8F"
$ws.Range("B74").Value = "#00CCCC"
$ws.Range("B75").Value = "#8C74BE"
$ws.Range("B76").Value = "#C31D76"
$ws.Range("B77").Value = "#363031"
$ws.Range("B78").Value = "#00B9B9"
$ws.Range("B79").Value = "#1FBCFF"
$ws.Range("B80").Value = "#E87377"
$ws.Range("B81").Value = "#F094C5"
$ws.Range("B82").Value = "#0081B8"
$ws.Range("C82").Value = "triangle"
$ws.Range("B83").Value = "#00ABF5"
$ws.Range("B84").Value = "#47346F"
$ws.Range("B85").Value = "#59418B"
$ws.Range("B86").Value = "#29B3E0"
$ws.Range("B87").Value = "#4CBFE6"
$ws.Range("B88").Value = "#5EC6E8"
$ws.Range("B89").Value = "#70CCEB"
$ws.Range("B90").Value = "#82D2ED"
$ws.Range("C90").Value = "circle"
$ws.Range("B91").Value = "#B89000"
$ws.Range("C91").Value = "circle"
$ws.Range("B92").Value = "#CCA000"
$ws.Range("C92").Value = "circle"
$ws.Range("B93").Value = "#FFDC5C"
$ws.Range("C93").Value = "circle"
$ws.Range("B94").Value = "#7B1417"
$ws.Range("C94").Value = "circle"
$ws.Range("B95").Value = "#95191C"
$ws.Range("C95").Value = "circle"
$ws.Range("B96").Value = "#AF1D21"
$ws.Range("C96").Value = "circle"
$ws.Range("B97").Value = "#C12025"
$ws.Range("C97").Value = "circle"
$ws.Range("B98").Value = "#201D1E"
$ws.Range("C98").Value = "circle"
$ws.Range("B99").Value = "#7A60A9"
$ws.Range("C99").Value = "circle"
$ws.Range("B100").Value = "#8A73B4"
$ws.Range("C100").Value = "circle"
$ws.Range("B101").Value = "#9B88BF"
$ws.Range("C101").Value = "circle"
$ws.Range("B102").Value = "#A695C6"
$ws.Range("C102").Value = "circle"
$ws.Range("B103").Value = "#B1A2CD"
$ws.Range("B104").Value = "#53817A"
$ws.Range("C104").Value = "diamond"
$ws.Range("B105").Value = "#00567A"
$ws.Range("C105").Value = "diamond"
$ws.Range("B106").Value = "#00648F"
$ws.Range("C106").Value = "diamond"
$ws.Range("B107").Value = "#0081B8"
$ws.Range("C107").Value = "diamond"
$ws.Range("B108").Value = "#E03E43"
$ws.Range("C108").Value = "diamond"
$ws.Range("B109").Value = "#009DE0"
$ws.Range("C109").Value = "diamond"
$ws.Range("B110").Value = "#6AA097"
$ws.Range("B111").Value = "#83AFA9"
$ws.Range("B112").Value = "#f57e20"
$ws.Range("C112").Value = "square"
$ws.Range("B113").Value = "#ACADB0"
$ws.Range("B114").Value = "#0AFF89"
$ws.Range("B115").Value = "#90DA61"
$ws.Range("C115").Value = "square"
$ws.Range("B116").Value = "#ADE48B"
$ws.Range("B117").Value = "#974a21"
$ws.Range("C117").Value = "square"
$ws.Range("B118").Value = "#7459B1"
$ws.Range("C118").Value = "circle"
$ws.Range("B119").Value = "#ffc808"
$ws.Range("C119").Value = "circle"
$ws.Range("B120").Value = "#D32228"
$ws.Range("C120").Value = "circle"
$ws.Range("B121").Value = "#88898D"
$ws.Range("C121").Value = "diamond"
$ws.Range("B122").Value = "#009DE0"
$ws.Range("C122").Value = "square"
$ws.Range("B123").Value = "#A01860"
$ws.Range("C123").Value = "triangle"
$ws.Range("B124").Value = "#00CCCC"
$ws.Range("C124").Value = "circle"
$ws.Range("B125").Value = "#7459B1"
$ws.Range("C125").Value = "triangle"
$ws.Range("B126").Value = "#363031"
$ws.Range("C126").Value = "triangle"
$ws.Range("B127").Value = "#00567A"
$ws.Range("C127").Value = "circle"

$ws.Range("E11").Select()